$d = $word.ActiveDocument

# 1. Add <w:noProof/> (NoProofing) to the runs that hold the three inline
#    drawings in the document.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.NoProofing = $true
    }
}

# 2. Locate the paragraph that ends with "...how do they want them named."
#    and append two new paragraphs of text right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($target -eq $null -and $p.Range.Text -like "*how do they want them named.*") {
        $target = $p
    }
}

$target.Range.InsertParagraphAfter()
$new1 = $target.Next()
$new1.Range.Text = "We’ll then need to repackage the files into a zip archive, rename it to a .xlsx file and upload it to the users device"

$new1.Range.InsertParagraphAfter()
$new2 = $new1.Next()
$new2.Range.Text = "It may be easier to build this as a desktop app, have the user’s PC pull the necessary data from the database and then parse and convert the files locally. Since there’s four phases to the process. Collecting the data from the data base, converting that data into the spreadsheet, zipping, and renaming the archive and finally downloading the file. The processing will cost the aquarium if done on a server, where downloading the data to the users machine and then doing all of the conversion locally will save money in that respect. "
